{"js": "// Split the single run of text in the Title, Author and Abstract\n// paragraphs into one run per word, with the in-between spaces each\n// becoming their own separate run (matching the target OOXML diff).\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\n// Build a <w:p> ooxml fragment that keeps the paragraph's pStyle (if any)\n// and represents `text` as alternating word/space runs.\nfunction buildParagraphOoxml(styleId, text) {\n  const words = text.split(\" \");\n  const runs = words\n    .map((w, i) => {\n      const wordRun = `<w:r><w:t xml:space=\"preserve\">${escapeXml(w)}</w:t></w:r>`;\n      if (i === words.length - 1) return wordRun;\n      const spaceRun = `<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>`;\n      return wordRun + spaceRun;\n    })\n    .join(\"\");\n  const pPr = styleId ? `<w:pPr><w:pStyle w:val=\"${escapeXml(styleId)}\"/></w:pPr>` : \"\";\n  return (\n    `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>` +\n    `<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">` +\n    `<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">` +\n    `<pkg:xmlData>` +\n    `<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">` +\n    `<w:body><w:p>${pPr}${runs}</w:p></w:body>` +\n    `</w:document>` +\n    `</pkg:xmlData></pkg:part></pkg:package>`\n  );\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// The three paragraphs (identified by style) whose single run of text\n// must be split into one run per word plus separate single-space runs.\nconst targetStyles = [\"Title\", \"Author\", \"Abstract\"];\n\nfor (const para of paragraphs.items) {\n  if (targetStyles.indexOf(para.style) !== -1) {\n    const originalText = para.text;\n    const ooxml = buildParagraphOoxml(para.style, originalText);\n    para.getRange().insertOoxml(ooxml, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Split the single run of text in the Title, Author and Abstract\n# paragraphs into one run per word, with the in-between spaces each\n# becoming their own separate run (matching the target OOXML diff).\n\n$d = $word.ActiveDocument\n\n$targetStyles = @(\"Title\", \"Author\", \"Abstract\")\n\nfunction Escape-Xml($s) {\n    return $s.Replace(\"&\", \"&amp;\").Replace(\"<\", \"&lt;\").Replace(\">\", \"&gt;\")\n}\n\nfunction Build-ParagraphXml($styleId, $text) {\n    $words = $text.Split(\" \")\n    $runs = \"\"\n    for ($wordIdx = 0; $wordIdx -lt $words.Count; $wordIdx++) {\n        $word = Escape-Xml($words[$wordIdx])\n        $runs = $runs + '<w:r><w:t xml:space=\"preserve\">' + $word + '</w:t></w:r>'\n        if ($wordIdx -lt ($words.Count - 1)) {\n            $runs = $runs + '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>'\n        }\n    }\n    $pPr = \"\"\n    if ($styleId) {\n        $pPr = '<w:pPr><w:pStyle w:val=\"' + (Escape-Xml($styleId)) + '\"/></w:pPr>'\n    }\n    $xml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + $pPr + $runs + '</w:p>'\n    return $xml\n}\n\n$count = $d.Paragraphs.Count\nfor ($pIdx = 1; $pIdx -le $count; $pIdx++) {\n    $p = $d.Paragraphs.Item($pIdx)\n    $styleName = $p.Style.NameLocal\n    if ($targetStyles -contains $styleName) {\n        # Range.Text includes the trailing paragraph mark (\\r); strip it\n        # before splitting into words.\n        $fullText = $p.Range.Text\n        $text = $fullText.Replace(\"`r\", \"\")\n        $xml = Build-ParagraphXml $styleName $text\n        [void]$p.Range.InsertXML($xml)\n    }\n}\n"}
